# Applies the commit "Remove CTI and DWI input, update the format and set
# default parameters to analyse test_data" to the parameters-analysis
# template workbook.
#
#   activity_analysis      (sheet 1)
#   connectivity_analysis  (sheet 2)
#   statistical_analysis   (sheet 3)

function RGB($r, $g, $b) { return ($b * 65536) + ($g * 256) + $r }

$wb = $excel.ActiveWorkbook

$wsActivity     = $wb.Worksheets.Item("activity_analysis")
$wsConnectivity = $wb.Worksheets.Item("connectivity_analysis")
$wsStats        = $wb.Worksheets.Item("statistical_analysis")

# ---------------------------------------------------------------------------
# 1. activity_analysis - switch the three "ers_erd" toggles (and the erp one)
#    from "on" back to "off" -- these become the new default parameters.
# ---------------------------------------------------------------------------
$wsActivity.Range("B2").Value  = "off"   # erp.sensor_enable
$wsActivity.Range("B9").Value  = "off"   # ers_erd.sensor_enable
$wsActivity.Range("B10").Value = "off"   # ers_erd.roi_enable
$wsActivity.Range("B12").Value = "off"   # ers_erd.mapping_enable

# Re-colour the existing on/off conditional formatting: red for "off",
# green for "on" (was a muted orange / grey pairing before).
$activityCF = $wsActivity.Range("B1:B1048576").FormatConditions
$activityOff = $activityCF.Item(1)
$activityOff.Font.Color = RGB(255, 0, 0)
$activityOn = $activityCF.Item(2)
$activityOn.Font.Color = RGB(0, 176, 80)

# ---------------------------------------------------------------------------
# 2. connectivity_analysis - turn the ICA connectivity (CTI/DWI-style) block
#    back on and consolidate / simplify the conditional formatting so a
#    single rule-pair covers the whole column instead of two separate
#    blocks (the single-cell override on B11 is removed).
# ---------------------------------------------------------------------------
$wsConnectivity.Range("B2").Value = "on"   # ica_conn.enable

$connFull = $wsConnectivity.Range("B1:B1048576")
$connFull.FormatConditions.Delete()

$connOff = $connFull.FormatConditions.Add(9, $null, $null, $null, "off", 2)
$connOff.Font.Color = RGB(255, 0, 0)

$connOn = $connFull.FormatConditions.Add(9, $null, $null, $null, "on", 2)
$connOn.Font.Color = RGB(0, 176, 80)

# ---------------------------------------------------------------------------
# 3. statistical_analysis - default to demeaning the individual maps, and
#    broaden the yes/no conditional formatting to the whole column.
# ---------------------------------------------------------------------------
$wsStats.Range("B4").Value = "yes"   # stats.demean

$statsFull = $wsStats.Range("B1:B1048576")
$statsFull.FormatConditions.Delete()

$statsNo = $statsFull.FormatConditions.Add(9, $null, $null, $null, "no", 0)
$statsNo.Font.Color = RGB(237, 125, 49)

$statsYes = $statsFull.FormatConditions.Add(9, $null, $null, $null, "yes", 0)
$statsYes.Font.Color = RGB(123, 123, 123)

# ---------------------------------------------------------------------------
# 4. Restore the per-sheet selections recorded in the workbook, activating
#    each sheet so the saved selection sticks, finishing on
#    statistical_analysis so it stays the active/selected tab.
# ---------------------------------------------------------------------------
$wsActivity.Activate() | Out-Null
$wsActivity.Range("B13").Select() | Out-Null

$wsConnectivity.Activate() | Out-Null
$wsConnectivity.Range("B3").Select() | Out-Null

$wsStats.Activate() | Out-Null
$wsStats.Range("B5").Select() | Out-Null
